# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 20:52"

# Swap "La Palma" (row 56) and "Lanzarote" (row 57), keeping B/C/D the same
# but toggling the Muertes (deaths) values between the two rows.
$ws.Range("A56").Value = "Lanzarote"
$ws.Range("E56").Value = 3

$ws.Range("A57").Value = "La Palma"
$ws.Range("E57").Value = 4
